$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 title change
$ws.Range("A1").Value = "cdwa1"

# Update numeric correlation values (B:F) and the .bag filename (G) for rows 4-23
$ws.Cells.Item(4, 2).Value = 0.5735415670576476
$ws.Cells.Item(4, 3).Value = 0.3892947992043204
$ws.Cells.Item(4, 4).Value = 0.4141080527342521
$ws.Cells.Item(4, 5).Value = 0.3219275759755983
$ws.Cells.Item(4, 6).Value = 0.6936566612013039
$ws.Cells.Item(4, 7).Value = "exp3_cdwa1_w3_C4_r0.bag"

$ws.Cells.Item(5, 2).Value = -0.4250783238122718
$ws.Cells.Item(5, 3).Value = -0.5524747170534537
$ws.Cells.Item(5, 4).Value = -0.5392094904281148
$ws.Cells.Item(5, 5).Value = 0.2737667704825002
$ws.Cells.Item(5, 6).Value = -0.3407558095979445
$ws.Cells.Item(5, 7).Value = "exp3_cdwa1_w3_C4_r1.bag"

$ws.Cells.Item(6, 2).Value = -0.3247574537317936
$ws.Cells.Item(6, 3).Value = -0.472956243854587
$ws.Cells.Item(6, 4).Value = -0.5597554112639272
$ws.Cells.Item(6, 5).Value = -0.5063208533582128
$ws.Cells.Item(6, 6).Value = -0.328693297293293
$ws.Cells.Item(6, 7).Value = "exp3_cdwa1_w3_C4_r2.bag"

$ws.Cells.Item(7, 2).Value = -0.1926057445151897
$ws.Cells.Item(7, 3).Value = -0.3852435535824921
$ws.Cells.Item(7, 4).Value = -0.4674057805024837
$ws.Cells.Item(7, 5).Value = -0.3893843817353697
$ws.Cells.Item(7, 6).Value = -0.2672051563932004
$ws.Cells.Item(7, 7).Value = "exp3_cdwa1_w3_C4_r3.bag"

$ws.Cells.Item(8, 2).Value = -0.7471526372682664
$ws.Cells.Item(8, 3).Value = -0.8581633615996517
$ws.Cells.Item(8, 4).Value = -0.8681092109053035
$ws.Cells.Item(8, 5).Value = -0.8425804111809766
$ws.Cells.Item(8, 6).Value = -0.8001141319985793
$ws.Cells.Item(8, 7).Value = "exp3_cdwa1_w3_C4_r4.bag"

$ws.Cells.Item(9, 2).Value = -0.665281799702523
$ws.Cells.Item(9, 3).Value = -0.8057258371340386
$ws.Cells.Item(9, 4).Value = -0.6562080555191648
$ws.Cells.Item(9, 5).Value = 0.2800614058751336
$ws.Cells.Item(9, 6).Value = -0.6732711784899386
$ws.Cells.Item(9, 7).Value = "exp3_cdwa1_w3_C8_r0.bag"

$ws.Cells.Item(10, 2).Value = -0.4848350264733896
$ws.Cells.Item(10, 3).Value = -0.6228208116816955
$ws.Cells.Item(10, 4).Value = -0.7314117991825992
$ws.Cells.Item(10, 5).Value = -0.6986377160198806
$ws.Cells.Item(10, 6).Value = -0.5118210612256765
$ws.Cells.Item(10, 7).Value = "exp3_cdwa1_w3_C8_r1.bag"

$ws.Cells.Item(11, 2).Value = -0.4705917243547721
$ws.Cells.Item(11, 3).Value = -0.6292998645976879
$ws.Cells.Item(11, 4).Value = -0.6853596076287582
$ws.Cells.Item(11, 5).Value = -0.6199312715133094
$ws.Cells.Item(11, 6).Value = -0.5279014313826627
$ws.Cells.Item(11, 7).Value = "exp3_cdwa1_w3_C8_r2.bag"

$ws.Cells.Item(12, 2).Value = -0.7985404213354906
$ws.Cells.Item(12, 3).Value = -0.8352846075245539
$ws.Cells.Item(12, 4).Value = -0.7623068358105568
$ws.Cells.Item(12, 5).Value = -0.4018847798969204
$ws.Cells.Item(12, 6).Value = -0.794006674163796
$ws.Cells.Item(12, 7).Value = "exp3_cdwa1_w3_C8_r3.bag"

$ws.Cells.Item(13, 2).Value = -0.931026362235881
$ws.Cells.Item(13, 3).Value = -0.9507707963141869
$ws.Cells.Item(13, 4).Value = -0.8761150428779862
$ws.Cells.Item(13, 5).Value = -0.7997713793771809
$ws.Cells.Item(13, 6).Value = -0.9296885435793761
$ws.Cells.Item(13, 7).Value = "exp3_cdwa1_w3_C8_r4.bag"

$ws.Cells.Item(14, 2).Value = -0.8035936605379675
$ws.Cells.Item(14, 3).Value = -0.6150892744134782
$ws.Cells.Item(14, 4).Value = -0.3467839242597464
$ws.Cells.Item(14, 5).Value = -0.1727501786359029
$ws.Cells.Item(14, 6).Value = -0.7587239428053292
$ws.Cells.Item(14, 7).Value = "exp3_cdwa1_w4_C4_r0.bag"

$ws.Cells.Item(15, 2).Value = -0.2982551266058114
$ws.Cells.Item(15, 3).Value = -0.6265152534246134
$ws.Cells.Item(15, 4).Value = -0.8234341956624245
$ws.Cells.Item(15, 5).Value = -0.8952555486569849
$ws.Cells.Item(15, 6).Value = -0.2067538804195689
$ws.Cells.Item(15, 7).Value = "exp3_cdwa1_w4_C4_r1.bag"

$ws.Cells.Item(16, 2).Value = -0.367638173200275
$ws.Cells.Item(16, 3).Value = -0.6113148475098945
$ws.Cells.Item(16, 4).Value = -0.6532333312306599
$ws.Cells.Item(16, 5).Value = -0.8607426421260935
$ws.Cells.Item(16, 6).Value = -0.4051541979841365
$ws.Cells.Item(16, 7).Value = "exp3_cdwa1_w4_C4_r2.bag"

$ws.Cells.Item(17, 2).Value = 0.2542959944535256
$ws.Cells.Item(17, 3).Value = 0.3900681680566253
$ws.Cells.Item(17, 4).Value = -0.597894267613539
$ws.Cells.Item(17, 5).Value = -0.6728806846632656
$ws.Cells.Item(17, 6).Value = 0.3679301498708003
$ws.Cells.Item(17, 7).Value = "exp3_cdwa1_w4_C4_r3.bag"

$ws.Cells.Item(18, 2).Value = -0.6560561650594285
$ws.Cells.Item(18, 3).Value = -0.7241762258201195
$ws.Cells.Item(18, 4).Value = -0.7952114182753782
$ws.Cells.Item(18, 5).Value = -0.7921233265666138
$ws.Cells.Item(18, 6).Value = -0.6872127552874909
$ws.Cells.Item(18, 7).Value = "exp3_cdwa1_w4_C4_r4.bag"

$ws.Cells.Item(19, 2).Value = -0.808200448917245
$ws.Cells.Item(19, 3).Value = -0.8991036442135437
$ws.Cells.Item(19, 4).Value = -0.8083167580656001
$ws.Cells.Item(19, 5).Value = -0.6287347071893497
$ws.Cells.Item(19, 6).Value = -0.8090974155071581
$ws.Cells.Item(19, 7).Value = "exp3_cdwa1_w4_C8_r0.bag"

$ws.Cells.Item(20, 2).Value = -0.4330522074341633
$ws.Cells.Item(20, 3).Value = -0.5805098481779875
$ws.Cells.Item(20, 4).Value = -0.693965518036944
$ws.Cells.Item(20, 5).Value = -0.7074790217815058
$ws.Cells.Item(20, 6).Value = -0.4304771318350327
$ws.Cells.Item(20, 7).Value = "exp3_cdwa1_w4_C8_r1.bag"

$ws.Cells.Item(21, 2).Value = -0.5674063411074832
$ws.Cells.Item(21, 3).Value = -0.7027551758585761
$ws.Cells.Item(21, 4).Value = -0.7049881677435962
$ws.Cells.Item(21, 5).Value = -0.8162605525629311
$ws.Cells.Item(21, 6).Value = -0.5096757580489818
$ws.Cells.Item(21, 7).Value = "exp3_cdwa1_w4_C8_r2.bag"

$ws.Cells.Item(22, 2).Value = -0.5088989650291121
$ws.Cells.Item(22, 3).Value = -0.7115117971173978
$ws.Cells.Item(22, 4).Value = -0.7830673043554186
$ws.Cells.Item(22, 5).Value = -0.8637539757257153
$ws.Cells.Item(22, 6).Value = -0.5032496601174382
$ws.Cells.Item(22, 7).Value = "exp3_cdwa1_w4_C8_r3.bag"

$ws.Cells.Item(23, 2).Value = -0.6849711479818195
$ws.Cells.Item(23, 3).Value = -0.7526202160098286
$ws.Cells.Item(23, 4).Value = -0.5749744164490402
$ws.Cells.Item(23, 5).Value = -0.7243958293798345
$ws.Cells.Item(23, 6).Value = -0.6473955811266633
$ws.Cells.Item(23, 7).Value = "exp3_cdwa1_w4_C8_r4.bag"

